$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume (E) cells to match the latest scrape.
# Cells whose new text happens to look like a plain number (e.g. "607.33")
# are forced to stay text cells (matching the original inlineStr/shared-string
# cell type) by temporarily applying a text number format, then restoring the
# original cell style so no visible formatting changes are introduced.

$rng = $ws.Range("D5")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '607.33'
$rng.Style = $origStyle

$rng = $ws.Range("D6")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '145.89'
$rng.Style = $origStyle

$rng = $ws.Range("D8")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.585'
$rng.Style = $origStyle

$rng = $ws.Range("D11")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '5.52'
$rng.Style = $origStyle

$rng = $ws.Range("D12")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.374'
$rng.Style = $origStyle

$rng = $ws.Range("D14")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '27.20'
$rng.Style = $origStyle

$rng = $ws.Range("D17")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.0000147'
$rng.Style = $origStyle

$rng = $ws.Range("D19")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '11.50'
$rng.Style = $origStyle

$rng = $ws.Range("D21")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '342.09'
$rng.Style = $origStyle

$rng = $ws.Range("D22")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '6.87'
$rng.Style = $origStyle

$rng = $ws.Range("D24")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '5.70'
$rng.Style = $origStyle

$rng = $ws.Range("D25")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '66.24'
$rng.Style = $origStyle

$rng = $ws.Range("D28")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '9.03'
$rng.Style = $origStyle

$rng = $ws.Range("D29")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '549.32'
$rng.Style = $origStyle

$rng = $ws.Range("D30")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.162'
$rng.Style = $origStyle

$rng = $ws.Range("D32")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '7.83'
$rng.Style = $origStyle

$rng = $ws.Range("D33")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '2.03'
$rng.Style = $origStyle

$rng = $ws.Range("D35")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '1.74'
$rng.Style = $origStyle

$rng = $ws.Range("D36")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '5.21'
$rng.Style = $origStyle

$rng = $ws.Range("D37")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '168.30'
$rng.Style = $origStyle

$rng = $ws.Range("D38")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '1.00'
$rng.Style = $origStyle

$rng = $ws.Range("D40")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '1.95'
$rng.Style = $origStyle

$rng = $ws.Range("D43")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '165.84'
$rng.Style = $origStyle

$rng = $ws.Range("D44")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '39.65'
$rng.Style = $origStyle

$rng = $ws.Range("D46")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '21.85'
$rng.Style = $origStyle

$rng = $ws.Range("D47")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.0563'
$rng.Style = $origStyle

$rng = $ws.Range("D48")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.624'
$rng.Style = $origStyle

$rng = $ws.Range("D49")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.0243'
$rng.Style = $origStyle

$rng = $ws.Range("D51")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '1.93'
$rng.Style = $origStyle

# Remaining cells (non-numeric-looking text, e.g. percentages or multi-dot prices)
# can be set directly without any special handling.

$ws.Range("D2").Value = '63.245.49'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '2.613.47'
$ws.Range("E3").Value = '  -2.01%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("E5").Value = '  +2.20%  '
$ws.Range("E6").Value = '  +1.44%  '
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("E8").Value = '  -0.40%  '
$ws.Range("D9").Value = '2.613.09'
$ws.Range("E9").Value = '  -1.61%  '
$ws.Range("E10").Value = '  +1.12%  '
$ws.Range("E11").Value = '  -3.01%  '
$ws.Range("E12").Value = '  +5.12%  '
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("E14").Value = '  -1.14%  '
$ws.Range("D15").Value = '3.088.56'
$ws.Range("E15").Value = '  -1.64%  '
$ws.Range("D16").Value = '63.123.35'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("E17").Value = '  +1.60%  '
$ws.Range("D18").Value = '2.612.04'
$ws.Range("E18").Value = '  -1.96%  '
$ws.Range("E19").Value = '  -0.25%  '
$ws.Range("E20").Value = '  +2.17%  '
$ws.Range("E21").Value = '  +0.72%  '
$ws.Range("E22").Value = '  +0.64%  '
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("E24").Value = '  -1.37%  '
$ws.Range("E25").Value = '  -1.77%  '
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("E27").Value = '  +3.81%  '
$ws.Range("E28").Value = '  +6.40%  '
$ws.Range("E29").Value = '  +1.84%  '
$ws.Range("E30").Value = '  -2.45%  '
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("E32").Value = '  -0.31%  '
$ws.Range("E33").Value = '  +2.24%  '
$ws.Range("D34").Value = '0.0₃0845'
$ws.Range("E34").Value = '  +3.86%  '
$ws.Range("E35").Value = '  -5.19%  '
$ws.Range("E36").Value = '  +1.94%  '
$ws.Range("E37").Value = '  -2.72%  '
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("E39").Value = '  -1.15%  '
$ws.Range("E40").Value = '  +6.08%  '
$ws.Range("E41").Value = '  -1.16%  '
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("E43").Value = '  -4.94%  '
$ws.Range("E44").Value = '  -1.36%  '
$ws.Range("E45").Value = '  -0.50%  '
$ws.Range("E46").Value = '  -2.01%  '
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("E48").Value = '  -1.79%  '
$ws.Range("E49").Value = '  +1.21%  '
$ws.Range("E50").Value = '  -0.91%  '
$ws.Range("E51").Value = '  +12.16%  '
